# Fix(PV) : centrage des pv
# Adds two new "Journal" rows (28 and 29) describing the "Placement" and
# "point de vie" tasks, with their duration (minutes) and explanation text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 28 - "Placement" task, 20 minutes
$ws.Range("A28").Value = "Placement"
$ws.Range("C28").Value = 20
$ws.Range("E28").Value = "Placment du joueur au debut"

# Row 29 - "point de vie" task, 10 minutes
$ws.Range("A29").Value = "point de vie"
$ws.Range("C29").Value = 10
$ws.Range("E29").Value = "Affichage de la vie "

# Move the view / active selection down to the newly edited row, like the
# author did while working on these rows.
$win = $excel.ActiveWindow
$win.ScrollRow = 22
$win.ScrollColumn = 5
$ws.Range("E29").Select()
